$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 780,
# pushing the former rows 780-842 down to 781-843 (dimension A1:R842 -> A1:R843).
$ws.Rows(780).Insert()

# Populate the newly inserted row 780 with the new record's data.
$ws.Cells.Item(780, 1).Value = 9
$ws.Cells.Item(780, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(780, 3).Value = "Metropolitana"
$ws.Cells.Item(780, 4).Value = 45265
$ws.Cells.Item(780, 5).Value = 13
$ws.Cells.Item(780, 6).Value = 100112024
$ws.Cells.Item(780, 7).Value = "Choclo"
$ws.Cells.Item(780, 8).Value = "Dulce o Americano"
$ws.Cells.Item(780, 9).Value = "Primera"
$ws.Cells.Item(780, 10).Value = 3400
$ws.Cells.Item(780, 11).Value = 600
$ws.Cells.Item(780, 12).Value = 700
$ws.Cells.Item(780, 13).Value = 650
$ws.Cells.Item(780, 14).Value = "`$/unidad"
$ws.Cells.Item(780, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(780, 16).Value = 650
$ws.Cells.Item(780, 17).Value = 1
$ws.Cells.Item(780, 18).Value = "Hortaliza"
